# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures in the
# cryptos list, as produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price values are plain decimal numbers (e.g. "151.28",
# "0.880"). Excel would normally auto-convert a bare numeric-looking string
# assigned through .Value into a real number (dropping significant trailing
# zeros / re-formatting exponents), so those particular cells are first
# switched to a Text number format to preserve the source text exactly.
$priceCellsToKeepAsText = @(
    "D5",
    "D6",
    "D13",
    "D14",
    "D20",
    "D21",
    "D22",
    "D24",
    "D28",
    "D29",
    "D31",
    "D33",
    "D34",
    "D35",
    "D41",
    "D42",
    "D45",
    "D46",
    "D47",
    "D49",
)
foreach ($addr in $priceCellsToKeepAsText) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "615.34"
$ws.Range("D6").Value = "151.28"
$ws.Range("D13").Value = "0.0000220"
$ws.Range("D14").Value = "32.06"
$ws.Range("D20").Value = "15.35"
$ws.Range("D21").Value = "444.75"
$ws.Range("D22").Value = "9.52"
$ws.Range("D24").Value = "77.36"
$ws.Range("D28").Value = "10.17"
$ws.Range("D29").Value = "8.49"
$ws.Range("D31").Value = "1.57"
$ws.Range("D33").Value = "0.165"
$ws.Range("D34").Value = "25.85"
$ws.Range("D35").Value = "6.13"
$ws.Range("D41").Value = "177.12"
$ws.Range("D42").Value = "0.0882"
$ws.Range("D45").Value = "0.880"
$ws.Range("D46").Value = "28.20"
$ws.Range("D47").Value = "44.98"
$ws.Range("D49").Value = "1.25"

# The remaining Price values already contain more than one "." so Excel
# cannot parse them as numbers and keeps them as text automatically.
$ws.Range("D2").Value = "67.445.79"
$ws.Range("D3").Value = "3.520.45"
$ws.Range("D7").Value = "3.519.20"
$ws.Range("D15").Value = "4.118.19"
$ws.Range("D16").Value = "3.521.72"
$ws.Range("D17").Value = "67.451.62"
$ws.Range("D26").Value = "3.660.97"
$ws.Range("D36").Value = "3.514.85"

# Volume(1h) (column E) values always remain text (they include the "%"
# sign and padding spaces so Excel never treats them as numbers).
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  +5.48%  "
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -2.85%  "
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("E31").Value = "  -6.71%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +3.91%  "
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("E44").Value = "  -3.84%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -4.13%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  -3.00%  "
